$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.757.46"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").Value = "3.477.12"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'567.37"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").Value = "'181.95"
$ws.Range("E6").Value = "  -3.67%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -3.22%  "

$ws.Range("D8").Value = "3.467.67"
$ws.Range("E8").Value = "  -2.88%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  +3.12%  "

$ws.Range("E11").Value = "  -3.29%  "

$ws.Range("D12").Value = "'53.57"
$ws.Range("E12").Value = "  -4.25%  "

$ws.Range("D13").Value = "'0.0000298"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "4.024.67"
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("D16").Value = "'19.09"
$ws.Range("E16").Value = "  -3.79%  "

$ws.Range("D17").Value = "68.587.44"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("D18").Value = "3.473.36"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").Value = "'536.60"
$ws.Range("E21").Value = "  +12.88%  "

$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("D23").Value = "'19.17"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").Value = "'4.96"
$ws.Range("E24").Value = "  -2.14%  "

$ws.Range("D25").Value = "'4.36"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").Value = "'93.93"
$ws.Range("E26").Value = "  +0.57%  "

$ws.Range("D27").Value = "'2.88"
$ws.Range("E27").Value = "  -4.57%  "

$ws.Range("D28").Value = "'10.72"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -3.47%  "

$ws.Range("D30").Value = "'31.32"
$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("E31").Value = "  -7.47%  "

$ws.Range("D32").Value = "'12.48"
$ws.Range("E32").Value = "  +2.50%  "

$ws.Range("D33").Value = "'64.24"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("E34").Value = "  -5.46%  "

$ws.Range("D35").Value = "'568.09"
$ws.Range("E35").Value = "  -2.46%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'37.58"
$ws.Range("E37").Value = "  -3.69%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.393"
$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +5.42%  "

$ws.Range("E40").Value = "  -4.80%  "

$ws.Range("E41").Value = "  -3.94%  "

$ws.Range("D42").Value = "'0.132"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").Value = "'3.31"
$ws.Range("E43").Value = "  -4.77%  "

$ws.Range("D44").Value = "3.209.06"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("E45").Value = "  -4.57%  "

$ws.Range("D46").Value = "'3.44"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").Value = "'8.96"
$ws.Range("E48").Value = "  -5.75%  "

$ws.Range("D49").Value = "'0.134"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").Value = "'0.997"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").Value = "'136.95"
$ws.Range("E51").Value = "  +0.36%  "
